$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1144.7646
$ws.Range("I70").Value = 1344.3334
$ws.Range("J70").Value = 920.25
$ws.Range("K70").Value = 4033.0002
$ws.Range("L70").Value = 2760.75
$ws.Range("M70").Value = -3763.0002
$ws.Range("N70").Value = -3300.75

$ws.Range("H73").Value = 1144.7646
$ws.Range("I73").Value = 1344.3334
$ws.Range("J73").Value = 920.25
$ws.Range("K73").Value = 4033.0002
$ws.Range("L73").Value = 2760.75
$ws.Range("M73").Value = -3097.0002
$ws.Range("N73").Value = -4632.75

$ws.Range("H111").Value = 6288.048
$ws.Range("I111").Value = 7142.7646
$ws.Range("J111").Value = 2655.5
$ws.Range("K111").Value = 21428.2938
$ws.Range("L111").Value = 7966.5
$ws.Range("M111").Value = -18361.2938
$ws.Range("N111").Value = -14100.5

$ws.Range("H116").Value = 2660
$ws.Range("I116").Value = 1980
$ws.Range("K116").Value = 1980
$ws.Range("M116").Value = 1462

$ws.Range("H129").Value = 2816.6667
$ws.Range("I129").Value = 6757.25
$ws.Range("J129").Value = 846.375
$ws.Range("K129").Value = 20271.75
$ws.Range("L129").Value = 2539.125
$ws.Range("M129").Value = -15271.75
$ws.Range("N129").Value = -12539.125

$ws.Range("H132").Value = 4241304
$ws.Range("I132").Value = 4549513.5
$ws.Range("J132").Value = 3425
$ws.Range("K132").Value = 13648540.5
$ws.Range("L132").Value = 10275
$ws.Range("M132").Value = -13646010.5
$ws.Range("N132").Value = -15335

$ws.Range("H138").Value = 1639.6154
$ws.Range("I138").Value = 1550.3889
$ws.Range("K138").Value = 4651.1667
$ws.Range("M138").Value = 488.8333000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2515.2964
$ws.Range("I61").Value = 2324.353
$ws.Range("J61").Value = 2839.9
$ws.Range("K61").Value = 2324.353
$ws.Range("L61").Value = 2839.9
$ws.Range("M61").Value = -2112.353
$ws.Range("N61").Value = -3263.9

$ws.Range("H97").Value = 29626.2
$ws.Range("I97").Value = 39151.117
$ws.Range("J97").Value = 2109.7778
$ws.Range("K97").Value = 39151.117
$ws.Range("L97").Value = 2109.7778
$ws.Range("M97").Value = -38655.117
$ws.Range("N97").Value = -3101.7778

$ws.Range("H102").Value = 102033.9
$ws.Range("I102").Value = 251794.75
$ws.Range("J102").Value = 2193.3333
$ws.Range("K102").Value = 251794.75
$ws.Range("L102").Value = 2193.3333
$ws.Range("M102").Value = -250172.75
$ws.Range("N102").Value = -5437.3333

$ws.Range("H125").Value = 48000
$ws.Range("J125").Value = 48000
$ws.Range("L125").Value = 48000
$ws.Range("N125").Value = -57840

$ws.Range("H136").Value = 2515.2964
$ws.Range("I136").Value = 2324.353
$ws.Range("J136").Value = 2839.9
$ws.Range("K136").Value = 6973.059
$ws.Range("L136").Value = 8519.700000000001
$ws.Range("M136").Value = -4423.059
$ws.Range("N136").Value = -13619.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 49981.383
$ws.Range("J20").Value = 1484.7142
$ws.Range("L20").Value = 1484.7142
$ws.Range("N20").Value = -1978.7142

$ws.Range("H99").Value = 1711.3
$ws.Range("I99").Value = 1181.7273
$ws.Range("J99").Value = 2017.8948
$ws.Range("K99").Value = 1181.7273
$ws.Range("L99").Value = 2017.8948
$ws.Range("M99").Value = 316.2727
$ws.Range("N99").Value = -5013.8948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3339.25
$ws.Range("I31").Value = 2033
$ws.Range("J31").Value = 4210.0835
$ws.Range("K31").Value = 2033
$ws.Range("L31").Value = 4210.0835
$ws.Range("M31").Value = -1738
$ws.Range("N31").Value = -4800.0835

$ws.Range("H34").Value = 3339.25
$ws.Range("I34").Value = 2033
$ws.Range("J34").Value = 4210.0835
$ws.Range("K34").Value = 2033
$ws.Range("L34").Value = 4210.0835
$ws.Range("M34").Value = -1831
$ws.Range("N34").Value = -4614.0835

$ws.Range("H44").Value = 37071
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 37071
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 37071
$ws.Range("N44").Value = -37955
$ws.Range("M44").ClearContents()

$ws.Range("H62").Value = 2621.2
$ws.Range("J62").Value = 2621.2
$ws.Range("L62").Value = 2621.2
$ws.Range("N62").Value = -3869.2

$ws.Range("H65").Value = 2621.2
$ws.Range("J65").Value = 2621.2
$ws.Range("L65").Value = 13106
$ws.Range("N65").Value = -19346

$ws.Range("H127").Value = 47500
$ws.Range("J127").Value = 47500
$ws.Range("L127").Value = 47500
$ws.Range("N127").Value = -57420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 761.98
$ws.Range("I131").Value = 340
$ws.Range("J131").Value = 784.18945
$ws.Range("K131").Value = 1020
$ws.Range("L131").Value = 2352.56835
$ws.Range("M131").Value = 4020
$ws.Range("N131").Value = -12432.56835

$ws.Range("H132").Value = 1706.1471
$ws.Range("I132").Value = 849.9
$ws.Range("J132").Value = 2062.9167
$ws.Range("K132").Value = 7649.099999999999
$ws.Range("L132").Value = 18566.2503
$ws.Range("M132").Value = -5119.099999999999
$ws.Range("N132").Value = -23626.2503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3266.25
$ws.Range("I126").Value = 4798.2
$ws.Range("J126").Value = 2172
$ws.Range("K126").Value = 14394.6
$ws.Range("L126").Value = 6516
$ws.Range("M126").Value = -11924.6
$ws.Range("N126").Value = -11456

$ws.Range("H132").Value = 3577.4666
$ws.Range("I132").Value = 3305.3333
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 9915.999899999999
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -7385.999899999999
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4117.1724
$ws.Range("I7").Value = 4523.8667
$ws.Range("J7").Value = 3681.4285
$ws.Range("K7").Value = 4523.8667
$ws.Range("L7").Value = 3681.4285
$ws.Range("M7").Value = -4411.8667
$ws.Range("N7").Value = -3905.4285

$ws.Range("H22").Value = 2671.5454
$ws.Range("I22").Value = 2597.5
$ws.Range("J22").Value = 2713.8572
$ws.Range("K22").Value = 2597.5
$ws.Range("L22").Value = 2713.8572
$ws.Range("M22").Value = -2302.5
$ws.Range("N22").Value = -3303.8572

$ws.Range("H27").Value = 2671.5454
$ws.Range("I27").Value = 2597.5
$ws.Range("J27").Value = 2713.8572
$ws.Range("K27").Value = 2597.5
$ws.Range("L27").Value = 2713.8572
$ws.Range("M27").Value = -2490.5
$ws.Range("N27").Value = -2927.8572

$ws.Range("H40").Value = 85583.836
$ws.Range("I40").Value = 251751.5
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 251751.5
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -251615.5
$ws.Range("N40").Value = -2772

$ws.Range("H55").Value = 860.9677
$ws.Range("I55").Value = 1382.5
$ws.Range("J55").Value = 531.5789
$ws.Range("K55").Value = 1382.5
$ws.Range("L55").Value = 531.5789
$ws.Range("M55").Value = -1209.5
$ws.Range("N55").Value = -877.5789

$ws.Range("H61").Value = 1924.9375
$ws.Range("I61").Value = 1957.1111
$ws.Range("J61").Value = 1883.5714
$ws.Range("K61").Value = 1957.1111
$ws.Range("L61").Value = 1883.5714
$ws.Range("M61").Value = -1755.1111
$ws.Range("N61").Value = -2287.5714

$ws.Range("H108").Value = 24500
$ws.Range("J108").Value = 24500
$ws.Range("L108").Value = 24500
$ws.Range("N108").Value = -32180

$ws.Range("H113").Value = 1924.9375
$ws.Range("I113").Value = 1957.1111
$ws.Range("J113").Value = 1883.5714
$ws.Range("K113").Value = 1957.1111
$ws.Range("L113").Value = 1883.5714
$ws.Range("M113").Value = 212.8888999999999
$ws.Range("N113").Value = -6223.5714

$ws.Range("H122").Value = 2449.7
$ws.Range("I122").Value = 2449.7
$ws.Range("K122").Value = 7349.099999999999
$ws.Range("M122").Value = -4899.099999999999

$ws.Range("H126").Value = 4117.1724
$ws.Range("I126").Value = 4523.8667
$ws.Range("J126").Value = 3681.4285
$ws.Range("K126").Value = 13571.6001
$ws.Range("L126").Value = 11044.2855
$ws.Range("M126").Value = -11101.6001
$ws.Range("N126").Value = -15984.2855

$ws.Range("H132").Value = 14850
$ws.Range("I132").Value = 16466.666
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 49399.99800000001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -46869.99800000001
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 167381.67
$ws.Range("I81").Value = 167278.17
$ws.Range("J81").Value = 167485.17
$ws.Range("K81").Value = 334556.34
$ws.Range("L81").Value = 334970.34
$ws.Range("M81").Value = -333495.34
$ws.Range("N81").Value = -337092.34

$ws.Range("H84").Value = 167381.67
$ws.Range("I84").Value = 167278.17
$ws.Range("J84").Value = 167485.17
$ws.Range("K84").Value = 1672781.7
$ws.Range("L84").Value = 1674851.7
$ws.Range("M84").Value = -1667477.7
$ws.Range("N84").Value = -1685459.7

$ws.Range("H113").Value = 709.1923
$ws.Range("I113").Value = 621.2778
$ws.Range("J113").Value = 907
$ws.Range("K113").Value = 1863.8334
$ws.Range("L113").Value = 2721
$ws.Range("M113").Value = 306.1666
$ws.Range("N113").Value = -7061

$ws.Range("H126").Value = 2955.6
$ws.Range("I126").Value = 3093.3333
$ws.Range("J126").Value = 2749
$ws.Range("K126").Value = 9279.999899999999
$ws.Range("L126").Value = 8247
$ws.Range("M126").Value = -6809.999899999999
$ws.Range("N126").Value = -13187
